$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:E) to (B:F)
$ws.Columns.Item(1).Insert()

# Copy the header formatting from B1 (original "A" header) onto the new A1 cell
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set new header text for column A
$ws.Range("A1").Value = "ID"

# Fill in the ID values for each data row
$ws.Range("A2").Value = "Hb 2"
$ws.Range("A3").Value = "Hb 3"
$ws.Range("A4").Value = "S 24"
$ws.Range("A5").Value = "S 28"
$ws.Range("A6").Value = "Hb 107"
$ws.Range("A7").Value = "Hb 66"
$ws.Range("A8").Value = "Hb 69"
$ws.Range("A9").Value = "Hb 95"
$ws.Range("A10").Value = "Hb 99"
$ws.Range("A11").Value = "Hb 92"
$ws.Range("A12").Value = "Hb 40"
$ws.Range("A13").Value = "Hb 41"
$ws.Range("A14").Value = "S 11"
$ws.Range("A15").Value = "Hb 57"
$ws.Range("A16").Value = "S 21"
$ws.Range("A17").Value = "S 22"
$ws.Range("A18").Value = "S 3"
$ws.Range("A19").Value = "S 4"
$ws.Range("A20").Value = "S 5"
$ws.Range("A21").Value = "Hb 74"
$ws.Range("A22").Value = "Hb 79"
$ws.Range("A23").Value = "Hb 32"
$ws.Range("A24").Value = "S 15"
$ws.Range("A25").Value = "S 16"

Write-Host "done"
